$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A23").Value = "fecha_DerAlc"
$ws.Range("A24").Value = "hora_DerAlc"
$ws.Range("C23").Value = "201:DERECHOS DE ALCOHOLEMIA"
$ws.Range("C24").Value = "201:DERECHOS DE ALCOHOLEMIA"

$ws.Range("D24").Select() | Out-Null
